# Add a new "Greece" market test-data sheet, built as a copy of the
# existing "Croatia" sheet (same layout/styles), with the market name and
# ticket reference swapped to the Greece values.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Select the whole sheet first so the source tab's leftover selection
# collapses to a full-sheet selection once it stops being the active tab
# (matches the artifact left behind by copying a sheet in the real UI).
$croatia.Cells.Select()

# Copy "Croatia" to a new tab placed immediately after it.
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Update the market name and ticket reference for the new sheet.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3189"

# Leave the active cell on B4, matching the new sheet's selection.
$greece.Range("B4").Select()
